$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for each data row (2-295).
# Update every row's value from 45189 (2023-09-20) to 45190 (2023-09-21).
$ws.Range("C2:C295").Value = 45190
